$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8 and 9 (the "candidates data" / "test results" online-test rows),
# which shifts the remaining rows (old 10-12) up to become rows 8-10.
$ws.Rows("8:9").Delete()

$ws.Range("D10").Select()
